$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 72
$ws.Cells.Item(72, 2).Value = 6811729
$ws.Cells.Item(72, 6).Value = "Stade Nyonnais"
$ws.Cells.Item(72, 7).Value = "FC Thun"
$ws.Cells.Item(72, 8).Value = 3
$ws.Cells.Item(72, 9).Value = 2
$ws.Cells.Item(72, 10).Value = "H"
$ws.Cells.Item(72, 11).Value = 3.4
$ws.Cells.Item(72, 12).Value = 3.5
$ws.Cells.Item(72, 13).Value = 2
$ws.Cells.Item(72, 14).Value = 3.8
$ws.Cells.Item(72, 15).Value = 3.75
$ws.Cells.Item(72, 16).Value = 1.85
$ws.Cells.Item(72, 17).Value = 0.5
$ws.Cells.Item(72, 18).Value = 1.975
$ws.Cells.Item(72, 19).Value = 1.825
$ws.Cells.Item(72, 20).Value = 3
$ws.Cells.Item(72, 21).Value = 1.9
$ws.Cells.Item(72, 22).Value = 1.9
$ws.Cells.Item(72, 23).Value = 2.8
$ws.Cells.Item(72, 24).Value = -1
$ws.Cells.Item(72, 25).Value = -1
$ws.Cells.Item(72, 26).Value = 0.9750000000000001
$ws.Cells.Item(72, 27).Value = -1
$ws.Cells.Item(72, 28).Value = 0.8999999999999999
$ws.Cells.Item(72, 29).Value = -1

# Row 73
$ws.Cells.Item(73, 2).Value = 6811942
$ws.Cells.Item(73, 6).Value = "FC Schaffhausen"
$ws.Cells.Item(73, 7).Value = "Neuchatel Xamax"
$ws.Cells.Item(73, 8).Value = 2
$ws.Cells.Item(73, 9).Value = 2
$ws.Cells.Item(73, 10).Value = "D"
$ws.Cells.Item(73, 11).Value = 3.5
$ws.Cells.Item(73, 12).Value = 3.5
$ws.Cells.Item(73, 13).Value = 1.95
$ws.Cells.Item(73, 14).Value = 4.2
$ws.Cells.Item(73, 15).Value = 3.3
$ws.Cells.Item(73, 16).Value = 1.95
$ws.Cells.Item(73, 17).Value = 0.5
$ws.Cells.Item(73, 18).Value = 1.85
$ws.Cells.Item(73, 19).Value = 1.95
$ws.Cells.Item(73, 20).Value = 2.5
$ws.Cells.Item(73, 21).Value = 1.925
$ws.Cells.Item(73, 22).Value = 1.875
$ws.Cells.Item(73, 23).Value = -1
$ws.Cells.Item(73, 24).Value = 2.3
$ws.Cells.Item(73, 25).Value = -1
$ws.Cells.Item(73, 26).Value = 0.8500000000000001
$ws.Cells.Item(73, 27).Value = -1
$ws.Cells.Item(73, 28).Value = 0.925
$ws.Cells.Item(73, 29).Value = -1

# Row 108
$ws.Cells.Item(108, 2).Value = 7617772
$ws.Cells.Item(108, 6).Value = "FC Vaduz"
$ws.Cells.Item(108, 7).Value = "FC Sion"
$ws.Cells.Item(108, 8).Value = 1
$ws.Cells.Item(108, 9).Value = 2
$ws.Cells.Item(108, 10).Value = "A"
$ws.Cells.Item(108, 11).Value = 3.75
$ws.Cells.Item(108, 12).Value = 3.6
$ws.Cells.Item(108, 13).Value = 1.833
$ws.Cells.Item(108, 14).Value = 5.25
$ws.Cells.Item(108, 15).Value = 4
$ws.Cells.Item(108, 16).Value = 1.615
$ws.Cells.Item(108, 17).Value = 0.75
$ws.Cells.Item(108, 18).Value = 2
$ws.Cells.Item(108, 19).Value = 1.8
$ws.Cells.Item(108, 20).Value = 2.75
$ws.Cells.Item(108, 21).Value = 1.825
$ws.Cells.Item(108, 22).Value = 1.975
$ws.Cells.Item(108, 23).Value = -1
$ws.Cells.Item(108, 24).Value = -1
$ws.Cells.Item(108, 25).Value = 0.615
$ws.Cells.Item(108, 26).Value = -0.5
$ws.Cells.Item(108, 27).Value = 0.4
$ws.Cells.Item(108, 28).Value = 0.4125
$ws.Cells.Item(108, 29).Value = -0.5

# Row 109
$ws.Cells.Item(109, 2).Value = 7617773
$ws.Cells.Item(109, 6).Value = "FC Thun"
$ws.Cells.Item(109, 7).Value = "Aarau"
$ws.Cells.Item(109, 8).Value = 1
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = "H"
$ws.Cells.Item(109, 11).Value = 1.727
$ws.Cells.Item(109, 12).Value = 3.8
$ws.Cells.Item(109, 13).Value = 4
$ws.Cells.Item(109, 14).Value = 1.7
$ws.Cells.Item(109, 15).Value = 4.2
$ws.Cells.Item(109, 16).Value = 4.5
$ws.Cells.Item(109, 17).Value = -0.75
$ws.Cells.Item(109, 18).Value = 1.85
$ws.Cells.Item(109, 19).Value = 1.95
$ws.Cells.Item(109, 20).Value = 3
$ws.Cells.Item(109, 21).Value = 1.9
$ws.Cells.Item(109, 22).Value = 1.9
$ws.Cells.Item(109, 23).Value = 0.7
$ws.Cells.Item(109, 24).Value = -1
$ws.Cells.Item(109, 25).Value = -1
$ws.Cells.Item(109, 26).Value = 0.425
$ws.Cells.Item(109, 27).Value = -0.5
$ws.Cells.Item(109, 28).Value = -1
$ws.Cells.Item(109, 29).Value = 0.8999999999999999

# Row 110
$ws.Cells.Item(110, 2).Value = 7617832
$ws.Cells.Item(110, 6).Value = "FC Baden"
$ws.Cells.Item(110, 7).Value = "FC Schaffhausen"
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 9).Value = 1
$ws.Cells.Item(110, 10).Value = "A"
$ws.Cells.Item(110, 11).Value = 2.9
$ws.Cells.Item(110, 12).Value = 3.4
$ws.Cells.Item(110, 13).Value = 2.2
$ws.Cells.Item(110, 14).Value = 2.9
$ws.Cells.Item(110, 15).Value = 3.5
$ws.Cells.Item(110, 16).Value = 2.375
$ws.Cells.Item(110, 17).Value = 0.25
$ws.Cells.Item(110, 18).Value = 1.775
$ws.Cells.Item(110, 19).Value = 2.025
$ws.Cells.Item(110, 20).Value = 2.75
$ws.Cells.Item(110, 21).Value = 1.95
$ws.Cells.Item(110, 22).Value = 1.85
$ws.Cells.Item(110, 23).Value = -1
$ws.Cells.Item(110, 24).Value = -1
$ws.Cells.Item(110, 25).Value = 1.375
$ws.Cells.Item(110, 26).Value = -1
$ws.Cells.Item(110, 27).Value = 1.025
$ws.Cells.Item(110, 28).Value = -1
$ws.Cells.Item(110, 29).Value = 0.8500000000000001

# Row 111
$ws.Cells.Item(111, 2).Value = 7617774
$ws.Cells.Item(111, 6).Value = "AC Bellinzona"
$ws.Cells.Item(111, 7).Value = "Wil 1900"
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = "D"
$ws.Cells.Item(111, 11).Value = 2.45
$ws.Cells.Item(111, 12).Value = 3.2
$ws.Cells.Item(111, 13).Value = 2.7
$ws.Cells.Item(111, 14).Value = 2.6
$ws.Cells.Item(111, 15).Value = 3.1
$ws.Cells.Item(111, 16).Value = 3
$ws.Cells.Item(111, 17).Value = 0
$ws.Cells.Item(111, 18).Value = 1.775
$ws.Cells.Item(111, 19).Value = 2.025
$ws.Cells.Item(111, 20).Value = 2
$ws.Cells.Item(111, 21).Value = 1.775
$ws.Cells.Item(111, 22).Value = 2.025
$ws.Cells.Item(111, 23).Value = -1
$ws.Cells.Item(111, 24).Value = 2.1
$ws.Cells.Item(111, 25).Value = -1
$ws.Cells.Item(111, 26).Value = 0
$ws.Cells.Item(111, 27).Value = 0
$ws.Cells.Item(111, 28).Value = -1
$ws.Cells.Item(111, 29).Value = 1.025

# Row 138
$ws.Cells.Item(138, 15).Value = 4
$ws.Cells.Item(138, 16).Value = 4
$ws.Cells.Item(138, 17).Value = -0.75
$ws.Cells.Item(138, 18).Value = 2
$ws.Cells.Item(138, 19).Value = 1.8

# Row 139
$ws.Cells.Item(139, 14).Value = 2.1
$ws.Cells.Item(139, 15).Value = 3.4
$ws.Cells.Item(139, 16).Value = 3.6
$ws.Cells.Item(139, 18).Value = 1.775
$ws.Cells.Item(139, 19).Value = 2.025
$ws.Cells.Item(139, 21).Value = 2
$ws.Cells.Item(139, 22).Value = 1.8

# Row 140
$ws.Cells.Item(140, 14).Value = 2.05
$ws.Cells.Item(140, 15).Value = 3.6
$ws.Cells.Item(140, 16).Value = 3.5
$ws.Cells.Item(140, 17).Value = -0.5
$ws.Cells.Item(140, 18).Value = 2
$ws.Cells.Item(140, 19).Value = 1.8

# Row 141
$ws.Cells.Item(141, 15).Value = 3.6
$ws.Cells.Item(141, 16).Value = 3.1
$ws.Cells.Item(141, 21).Value = 1.85
$ws.Cells.Item(141, 22).Value = 1.95

# Row 142
$ws.Cells.Item(142, 16).Value = 6.5
$ws.Cells.Item(142, 18).Value = 1.775
$ws.Cells.Item(142, 19).Value = 2.025
$ws.Cells.Item(142, 21).Value = 1.975
$ws.Cells.Item(142, 22).Value = 1.825
